# Apply edit: "upload project with multiple files per sample"
# For the TB123 sample row (row 3), the BAM and FASTQ file columns are
# cleared to "NULL" since this sample no longer references single files
# directly (multiple files per sample are now tracked elsewhere).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "NULL"
$ws.Range("C3").Value = "NULL"
$ws.Range("D3").Value = "NULL"

# Update the active selection to reflect the last-edited cell.
$ws.Range("D3").Select()
